$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.032.83"
$ws.Range("E2").Value = "  -2.00%  "
$ws.Range("D3").Value = "3.497.29"
$ws.Range("E3").Value = "  -2.06%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'579.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.56%  "
$ws.Range("D6").Value = "'177.57"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.05%  "
$ws.Range("E7").Value = "  +3.90%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "'0.636"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.18%  "
$ws.Range("D10").Value = "'0.161"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.16%  "
$ws.Range("D11").Value = "'55.80"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.37%  "
$ws.Range("E12").Value = "  +2.28%  "
$ws.Range("D13").Value = "'9.26"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.74%  "
$ws.Range("D14").Value = "4.052.94"
$ws.Range("E14").Value = "  -1.42%  "
$ws.Range("D15").Value = "3.503.65"
$ws.Range("E15").Value = "  -1.23%  "
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("D17").Value = "'18.27"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.11%  "
$ws.Range("D18").Value = "65.903.62"
$ws.Range("E18").Value = "  -2.19%  "
$ws.Range("E19").Value = "  -0.05%  "
$ws.Range("E20").Value = "  +1.34%  "
$ws.Range("D21").Value = "'410.20"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.40%  "
$ws.Range("D22").Value = "'4.25"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +7.97%  "
$ws.Range("D23").Value = "'4.50"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +8.93%  "
$ws.Range("D24").Value = "'84.86"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.77%  "
$ws.Range("D25").Value = "'13.31"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +9.90%  "
$ws.Range("D26").Value = "'11.05"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.65%  "
$ws.Range("E27").Value = "  -2.08%  "
$ws.Range("D28").Value = "'6.03"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.88%  "
$ws.Range("D29").Value = "'9.14"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.31%  "
$ws.Range("D30").Value = "'30.27"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.75%  "
$ws.Range("E31").Value = "  -0.15%  "
$ws.Range("E32").Value = "  -0.42%  "
$ws.Range("D33").Value = "'592.46"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -8.20%  "
$ws.Range("E34").Value = "  -1.51%  "
$ws.Range("D35").Value = "'60.92"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.78%  "
$ws.Range("E36").Value = "  +1.87%  "
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").Value = "'3.67"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +9.66%  "
$ws.Range("D39").Value = "0.0₃0795"
$ws.Range("E39").Value = "  -4.03%  "
$ws.Range("E40").Value = "  -5.11%  "
$ws.Range("E41").Value = "  -2.31%  "
$ws.Range("D42").Value = "3.229.34"
$ws.Range("E42").Value = "  +5.70%  "
$ws.Range("E43").Value = "  +0.19%  "
$ws.Range("E44").Value = "  +2.39%  "
$ws.Range("D45").Value = "'3.32"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.49%  "
$ws.Range("D46").Value = "'2.55"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.96%  "
$ws.Range("D47").Value = "'0.0419"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.11%  "
$ws.Range("E48").Value = "  -6.04%  "
$ws.Range("E49").Value = "  +1.01%  "
$ws.Range("D50").Value = "'8.54"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.47%  "
$ws.Range("D51").Value = "'137.64"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.64%  "
